$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.902.80"
$ws.Range("E2").Value = "  -1.60%  "
$ws.Range("D3").Value = "1.812.99"
$ws.Range("E3").Value = "  -0.91%  "
$ws.Range("D4").Value = "'1.003"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.22%  "
$ws.Range("D5").Value = "'310.59"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.38%  "
$ws.Range("D6").Value = "'1.003"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.18%  "
$ws.Range("D7").Value = "'0.4481"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +4.71%  "
$ws.Range("D8").Value = "'0.3679"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.60%  "
$ws.Range("D9").Value = "'0.07279"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.28%  "
$ws.Range("D10").Value = "'0.8523"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.72%  "
$ws.Range("D11").Value = "'20.64"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.32%  "
$ws.Range("D12").Value = "1.797.80"
$ws.Range("E12").Value = "  -1.63%  "
$ws.Range("D13").Value = "'6.608"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.70%  "
$ws.Range("D14").Value = "'0.07105"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.12%  "
$ws.Range("D15").Value = "'5.303"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.39%  "
$ws.Range("D16").Value = "'91.71"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +3.03%  "
$ws.Range("D17").Value = "'1.005"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.17%  "
$ws.Range("D18").Value = "'0.000008740"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.51%  "
$ws.Range("E19").Value = "  -0.34%  "
$ws.Range("D20").Value = "'14.89"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.41%  "
$ws.Range("D21").Value = "26.929.70"
$ws.Range("E21").Value = "  -1.40%  "
$ws.Range("D22").Value = "'5.148"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.00%  "
$ws.Range("D23").Value = "'10.89"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.00%  "
$ws.Range("D24").Value = "'1.986"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.17%  "
$ws.Range("D25").Value = "'151.76"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.71%  "
$ws.Range("D26").Value = "'2.219"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.02%  "
$ws.Range("D27").Value = "'18.41"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.03%  "
$ws.Range("D28").Value = "'5.211"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.93%  "
$ws.Range("D29").Value = "'116.23"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.52%  "
$ws.Range("D30").Value = "'0.08831"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.80%  "
$ws.Range("D31").Value = "'1.173"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.64%  "
$ws.Range("D32").Value = "'0.7465"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.90%  "
$ws.Range("D33").Value = "'2.951"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +4.58%  "
$ws.Range("D34").Value = "'4.426"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.07%  "
$ws.Range("D35").Value = "'1.002"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.26%  "
$ws.Range("D36").Value = "'1.088"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -2.28%  "
$ws.Range("D37").Value = "'0.01954"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.35%  "
$ws.Range("D38").Value = "'0.05202"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.65%  "
$ws.Range("D39").Value = "'0.5276"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +3.71%  "
$ws.Range("D40").Value = "'2.866"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.71%  "
$ws.Range("D41").Value = "'7.089"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.52%  "
$ws.Range("D42").Value = "'0.1692"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.58%  "
$ws.Range("D43").Value = "'0.5199"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +8.78%  "
$ws.Range("D44").Value = "'8.413"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -3.41%  "
$ws.Range("D45").Value = "'10.50"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.75%  "
$ws.Range("E46").Value = "  +6.00%  "
$ws.Range("D47").Value = "'105.08"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.60%  "
$ws.Range("D48").Value = "'1.002"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.26%  "
$ws.Range("B49").Value = "Cronos"
$ws.Range("C49").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D49").Value = "'0.06375"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.12%  "
$ws.Range("B50").Value = "NEARProtocol"
$ws.Range("C50").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D50").Value = "'1.655"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.86%  "
$ws.Range("D51").Value = "'0.9155"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.64%  "
